$d = $word.ActiveDocument

# Replace protocol number
$d.Content.Find.Execute("307/2023", $false, $false, $false, $false, $false, $true, 1, $false, "223/2025", 2)

Write-Host "done"
